$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before row 1063 (shifts old 1063-1130 down to 1068-1135)
$ws.Range("A1063:A1067").EntireRow.Insert()

# New data rows (1063-1067): Tuna prices for 44610, Región Metropolitana
$data = @(
    @(1063, 44610, "Tuna", "Extra",   1900, 700, 750, 729),
    @(1064, 44610, "Tuna", "Primera", 4300, 500, 600, 542),
    @(1065, 44610, "Tuna", "Segunda", 3300, 400, 450, 421),
    @(1066, 44610, "Tuna", "Super",   2800, 800, 900, 846),
    @(1067, 44610, "Tuna", "Tercera", 3400, 200, 300, 232)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = 6
    $ws.Cells.Item($r, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($r, 3).Value = "Metropolitana"
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = 13
    $ws.Cells.Item($r, 6).Value = 100112027
    $ws.Cells.Item($r, 7).Value = "Melón"
    $ws.Cells.Item($r, 8).Value = $row[2]
    $ws.Cells.Item($r, 9).Value = $row[3]
    $ws.Cells.Item($r, 10).Value = $row[4]
    $ws.Cells.Item($r, 11).Value = $row[5]
    $ws.Cells.Item($r, 12).Value = $row[6]
    $ws.Cells.Item($r, 13).Value = $row[7]
    $ws.Cells.Item($r, 14).Value = "`$/unidad"
    $ws.Cells.Item($r, 15).Value = "Región Metropolitana"
    $ws.Cells.Item($r, 16).Value = $row[7]
    $ws.Cells.Item($r, 17).Value = 1
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}

